$wb = $excel.ActiveWorkbook

# 1. Rename "Sheet1" to "Comparison" (the defined name
#    _xlnm._FilterDatabase that points at this sheet follows the rename
#    automatically since it is stored as a sheet-relative reference).
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Comparison"

# 2. Re-write the per-row LEN($E<row>) formulas in column D as shared
#    formulas, matching the grouping Excel produced when the formula was
#    authored via fill-down in chunks of 64 rows.
$ws.Range("D2:D65").Formula   = "=LEN(`$E2)"
$ws.Range("D66:D129").Formula = "=LEN(`$E66)"
$ws.Range("D130:D193").Formula = "=LEN(`$E130)"
$ws.Range("D194:D257").Formula = "=LEN(`$E194)"
$ws.Range("D258:D321").Formula = "=LEN(`$E258)"
$ws.Range("D322:D387").Formula = "=LEN(`$E322)"
